$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.159.10'
$ws.Range('E2').Value = '  +1.48%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.137.07'
$ws.Range('E3').Value = '  +1.37%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.81'
$ws.Range('E5').Value = '  +2.58%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.18'
$ws.Range('E6').Value = '  +1.95%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.506'
$ws.Range('E8').Value = '  +10.32%  '

$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('E10').Value = '  +2.31%  '

$ws.Range('E11').Value = '  +4.05%  '

$ws.Range('E12').Value = '  +3.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.676.30'
$ws.Range('E13').Value = '  +1.41%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.75'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000171'
$ws.Range('E15').Value = '  +5.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.208.55'
$ws.Range('E16').Value = '  +1.38%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.24'
$ws.Range('E17').Value = '  +5.67%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.142.30'
$ws.Range('E18').Value = '  +1.61%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.96'
$ws.Range('E19').Value = '  +3.69%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.18'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.17'
$ws.Range('E21').Value = '  +7.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.02%  '

$ws.Range('E23').Value = '  -0.79%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.35'
$ws.Range('E24').Value = '  +2.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.513'
$ws.Range('E25').Value = '  +2.74%  '

$ws.Range('E26').Value = '  +0.41%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.30%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0888'
$ws.Range('E28').Value = '  +2.27%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.80'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.22'
$ws.Range('E30').Value = '  +5.70%  '

$ws.Range('E31').Value = '  +0.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.78'
$ws.Range('E32').Value = '  +3.96%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.17'
$ws.Range('E33').Value = '  +6.04%  '

$ws.Range('E34').Value = '  +3.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.02'
$ws.Range('E35').Value = '  +1.16%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.24'
$ws.Range('E36').Value = '  +3.45%  '

$ws.Range('E37').Value = '  +8.69%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.68'
$ws.Range('E38').Value = '  -0.19%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.68'
$ws.Range('E39').Value = '  +6.11%  '

$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.615.85'
$ws.Range('E40').Value = '  +9.24%  '

$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0676'
$ws.Range('E41').Value = '  +2.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.22'
$ws.Range('E42').Value = '  +3.73%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.82'
$ws.Range('E43').Value = '  +5.74%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.700'
$ws.Range('E44').Value = '  +0.81%  '

$ws.Range('E45').Value = '  +2.54%  '

$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.23'
$ws.Range('E47').Value = '  +4.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.983'
$ws.Range('E48').Value = '  +2.55%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0990'
$ws.Range('E49').Value = '  +9.96%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.25'
$ws.Range('E50').Value = '  +2.75%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.751'
$ws.Range('E51').Value = '  -1.55%  '

Write-Output "Updates applied"